$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 59. This pushes the existing rows
# 59-67 down to 61-69, preserving all of their data/formatting.
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

# Fill the two newly inserted (blank) rows with the new weekly data.

# New row 59: Murcott, Primera
$ws.Range("A59").Value = 11
$ws.Range("B59").Value = "Vega Monumental Concepción"
$ws.Range("C59").Value = "Bíobío"
$ws.Range("D59").Value = 44461
$ws.Range("E59").Value = 8
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100102
$ws.Range("H59").Value = "Cítricos"
$ws.Range("I59").Value = 100102004
$ws.Range("J59").Value = "Mandarina"
$ws.Range("K59").Value = "Murcott"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 100
$ws.Range("N59").Value = 6000
$ws.Range("O59").Value = 7000
$ws.Range("P59").Value = 6500
$ws.Range("Q59").Value = "$/bandeja 10 kilos"
$ws.Range("R59").Value = "Provincia de Limarí"
$ws.Range("S59").Value = 650
$ws.Range("T59").Value = 10

# New row 60: Murcott, Segunda
$ws.Range("A60").Value = 11
$ws.Range("B60").Value = "Vega Monumental Concepción"
$ws.Range("C60").Value = "Bíobío"
$ws.Range("D60").Value = 44461
$ws.Range("E60").Value = 8
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100102
$ws.Range("H60").Value = "Cítricos"
$ws.Range("I60").Value = 100102004
$ws.Range("J60").Value = "Mandarina"
$ws.Range("K60").Value = "Murcott"
$ws.Range("L60").Value = "Segunda"
$ws.Range("M60").Value = 50
$ws.Range("N60").Value = 5000
$ws.Range("O60").Value = 5000
$ws.Range("P60").Value = 5000
$ws.Range("Q60").Value = "$/bandeja 10 kilos"
$ws.Range("R60").Value = "Provincia de Limarí"
$ws.Range("S60").Value = 500
$ws.Range("T60").Value = 10
